$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 22422200
$ws.Range("E8").Value = 16716400
$ws.Range("F8").Value = 18958900
$ws.Range("G8").Value = 18473200
$ws.Range("H8").Value = 19707200
$ws.Range("I8").Value = 19782300
$ws.Range("J8").Value = 19729100

$ws.Range("D9").Value = 14970400
$ws.Range("E9").Value = 9502300
$ws.Range("F9").Value = 9999900
$ws.Range("G9").Value = 10915200
$ws.Range("H9").Value = 12261700
$ws.Range("I9").Value = 13698800
$ws.Range("J9").Value = 15091600

$ws.Range("D10").Value = 7451800
$ws.Range("E10").Value = 7214100
$ws.Range("F10").Value = 8959000
$ws.Range("G10").Value = 7558000
$ws.Range("H10").Value = 7445500
$ws.Range("I10").Value = 6083400
$ws.Range("J10").Value = 4637600

$ws.Range("D14").Value = -7800
$ws.Range("J14").Value = 55300

$ws.Range("D15").Value = 2995000
$ws.Range("E15").Value = 2198800
$ws.Range("F15").Value = 2138800
$ws.Range("G15").Value = 1728500
$ws.Range("H15").Value = 1676100
$ws.Range("I15").Value = 1637400
$ws.Range("J15").Value = 1761100

$ws.Range("D17").Value = 21051600
$ws.Range("E17").Value = 13988900
$ws.Range("F17").Value = 14633800
$ws.Range("G17").Value = 14722200
$ws.Range("H17").Value = 16128900
$ws.Range("I17").Value = 17265700
$ws.Range("J17").Value = 18430900

$ws.Range("D18").Value = 1370700
$ws.Range("E18").Value = 2727400
$ws.Range("F18").Value = 4325100
$ws.Range("G18").Value = 3751000
$ws.Range("H18").Value = 3578300
$ws.Range("I18").Value = 2516600
$ws.Range("J18").Value = 1298200

$ws.Range("D20").Value = 492000
$ws.Range("E20").Value = 334400
$ws.Range("F20").Value = 261300
$ws.Range("G20").Value = 235800
$ws.Range("H20").Value = 163100
$ws.Range("I20").Value = 121200
$ws.Range("J20").Value = 154200

$ws.Range("D21").Value = 4879900
$ws.Range("E21").Value = 5282500
$ws.Range("F21").Value = 6742800
$ws.Range("G21").Value = 5733200
$ws.Range("H21").Value = 5437200
$ws.Range("I21").Value = 4290900
$ws.Range("J21").Value = "NA"

$ws.Range("D22").Value = 1446800
$ws.Range("E22").Value = 1011800
$ws.Range("F22").Value = 1179200
$ws.Range("G22").Value = 1159700
$ws.Range("H22").Value = 1155700
$ws.Range("I22").Value = 1320400
$ws.Range("J22").Value = 1148100

$ws.Range("D23").Value = 415800
$ws.Range("E23").Value = 2050000
$ws.Range("F23").Value = 3407200
$ws.Range("G23").Value = 2827100
$ws.Range("H23").Value = 2585700
$ws.Range("I23").Value = 1317400
$ws.Range("J23").Value = 304300

$ws.Range("D24").Value = 180700
$ws.Range("E24").Value = 514300
$ws.Range("F24").Value = 845800
$ws.Range("G24").Value = 814400
$ws.Range("H24").Value = 671200
$ws.Range("I24").Value = 372600
$ws.Range("J24").Value = 129000

$ws.Range("D26").Value = 235100
$ws.Range("E26").Value = 1535700
$ws.Range("F26").Value = 2561400
$ws.Range("G26").Value = 2012800
$ws.Range("H26").Value = 1914500
$ws.Range("I26").Value = 944800
$ws.Range("J26").Value = 175300

$ws.Range("D27").Value = 224300
$ws.Range("E27").Value = 1264500
$ws.Range("F27").Value = 2026100
$ws.Range("G27").Value = 1596500
$ws.Range("H27").Value = 1547300
$ws.Range("I27").Value = 818100
$ws.Range("J27").Value = 175200

$ws.Range("D32").Value = -492000
$ws.Range("E32").Value = -334400
$ws.Range("F32").Value = -261300
$ws.Range("G32").Value = -235800
$ws.Range("H32").Value = -163100
$ws.Range("I32").Value = -121200
$ws.Range("J32").Value = -154200

$ws.Range("D33").Value = 224300
$ws.Range("E33").Value = 1264500
$ws.Range("F33").Value = 2026100
$ws.Range("G33").Value = 1596500
$ws.Range("H33").Value = 1547300
$ws.Range("I33").Value = 818100
$ws.Range("J33").Value = 175200

$ws.Range("D35").Value = 224300
$ws.Range("E35").Value = 1264500
$ws.Range("F35").Value = 2026100
$ws.Range("G35").Value = 1596500
$ws.Range("H35").Value = 1547300
$ws.Range("I35").Value = 818100
$ws.Range("J35").Value = 175200

$ws.Range("D41").Value = 1377600
$ws.Range("E41").Value = 1159200
$ws.Range("F41").Value = 1109800
$ws.Range("G41").Value = 1816300
$ws.Range("H41").Value = 1386400
$ws.Range("I41").Value = 1559100
$ws.Range("J41").Value = 1269300

$ws.Range("D42").Value = 91300
$ws.Range("I42").Value = 13900
$ws.Range("J42").Value = 14300

$ws.Range("D43").Value = 4481100
$ws.Range("E43").Value = 3001300
$ws.Range("F43").Value = 2846100
$ws.Range("G43").Value = 2711900
$ws.Range("H43").Value = 2478000
$ws.Range("I43").Value = 2515900
$ws.Range("J43").Value = 2764000

$ws.Range("D44").Value = 1096100
$ws.Range("E44").Value = 1020900
$ws.Range("F44").Value = 804800
$ws.Range("G44").Value = 994700
$ws.Range("H44").Value = 960100
$ws.Range("I44").Value = 1042200
$ws.Range("J44").Value = 1116900

$ws.Range("D45").Value = 157400
$ws.Range("E45").Value = 304800
$ws.Range("F45").Value = 220700
$ws.Range("G45").Value = 96700
$ws.Range("H45").Value = 249200
$ws.Range("I45").Value = 224500
$ws.Range("J45").Value = 240200

$ws.Range("D46").Value = 7203500
$ws.Range("E46").Value = 5486200
$ws.Range("F46").Value = 4981400
$ws.Range("G46").Value = 5619600
$ws.Range("H46").Value = 5073700
$ws.Range("I46").Value = 5355600
$ws.Range("J46").Value = 5404700

$ws.Range("D47").Value = 3669600
$ws.Range("E47").Value = 3747800
$ws.Range("F47").Value = 4006200
$ws.Range("G47").Value = 3482800
$ws.Range("H47").Value = 3104200
$ws.Range("I47").Value = 2619400
$ws.Range("J47").Value = 2358100

$ws.Range("D48").Value = 42229400
$ws.Range("E48").Value = 33104600
$ws.Range("F48").Value = 32601700
$ws.Range("G48").Value = 27957300
$ws.Range("H48").Value = 26923800
$ws.Range("I48").Value = 26270600
$ws.Range("J48").Value = 26412200

$ws.Range("D49").Value = 3227100
$ws.Range("E49").Value = 2677700
$ws.Range("F49").Value = 2590300
$ws.Range("G49").Value = 2644900
$ws.Range("H49").Value = 2817400
$ws.Range("I49").Value = 3730600
$ws.Range("J49").Value = 3626400

$ws.Range("D52").Value = 2528300
$ws.Range("E52").Value = 1709000
$ws.Range("F52").Value = 1659200
$ws.Range("G52").Value = 1133600
$ws.Range("H52").Value = 998900
$ws.Range("I52").Value = 476900
$ws.Range("J52").Value = 401600

$ws.Range("D54").Value = 58857800
$ws.Range("E54").Value = 46725400
$ws.Range("F54").Value = 45838900
$ws.Range("G54").Value = 40838200
$ws.Range("H54").Value = 38918000
$ws.Range("I54").Value = 38453100
$ws.Range("J54").Value = 38203100

$ws.Range("D57").Value = 2299800
$ws.Range("E57").Value = 1789700
$ws.Range("F57").Value = 1395500
$ws.Range("G57").Value = 3714300
$ws.Range("H57").Value = 1822200
$ws.Range("I57").Value = 1091400
$ws.Range("J57").Value = 1353900

$ws.Range("D58").Value = 16831900
$ws.Range("E58").Value = 14519800
$ws.Range("F58").Value = 13811500
$ws.Range("G58").Value = 11010100
$ws.Range("H58").Value = 10620200
$ws.Range("I58").Value = 10677900
$ws.Range("J58").Value = 10296300

$ws.Range("D59").Value = 4012800
$ws.Range("E59").Value = 3013000
$ws.Range("F59").Value = 3171600
$ws.Range("G59").Value = 835800
$ws.Range("H59").Value = 2247100
$ws.Range("I59").Value = 2121000
$ws.Range("J59").Value = 2685800

$ws.Range("D60").Value = 23144600
$ws.Range("E60").Value = 19322400
$ws.Range("F60").Value = 18378600
$ws.Range("G60").Value = 15560200
$ws.Range("H60").Value = 14689500
$ws.Range("I60").Value = 13890300
$ws.Range("J60").Value = 14336100

$ws.Range("D61").Value = 18425700
$ws.Range("E61").Value = 11530500
$ws.Range("F61").Value = 11604600
$ws.Range("G61").Value = 12161900
$ws.Range("H61").Value = 12519300
$ws.Range("I61").Value = 14165700
$ws.Range("J61").Value = 14499600

$ws.Range("D62").Value = 1316500
$ws.Range("E62").Value = 706900
$ws.Range("F62").Value = 763300
$ws.Range("G62").Value = 533500
$ws.Range("H62").Value = 549800
$ws.Range("I62").Value = 607900
$ws.Range("J62").Value = 528400

$ws.Range("D66").Value = 45851000
$ws.Range("E66").Value = 33961700
$ws.Range("F66").Value = 33351300
$ws.Range("G66").Value = 30430300
$ws.Range("H66").Value = 29649700
$ws.Range("I66").Value = 30122800
$ws.Range("J66").Value = 30651600

$ws.Range("D72").Value = 6685100
$ws.Range("E72").Value = 7115000
$ws.Range("F72").Value = 6910800
$ws.Range("G72").Value = 5698000
$ws.Range("H72").Value = 4887000
$ws.Range("I72").Value = 3671100
$ws.Range("J72").Value = 3006700

$ws.Range("D76").Value = 13006900
$ws.Range("E76").Value = 12763700
$ws.Range("F76").Value = 12487500
$ws.Range("G76").Value = 10407900
$ws.Range("H76").Value = 9268300
$ws.Range("I76").Value = 8330300
$ws.Range("J76").Value = 7551500

$ws.Range("D81").Value = 224300
$ws.Range("E81").Value = 1264500
$ws.Range("F81").Value = 2026100
$ws.Range("G81").Value = 1596500
$ws.Range("H81").Value = 1547300
$ws.Range("I81").Value = 818100
$ws.Range("J81").Value = 175200

$ws.Range("D83").Value = 3011900
$ws.Range("E83").Value = 2216800
$ws.Range("F83").Value = 2152600
$ws.Range("G83").Value = 1743300
$ws.Range("H83").Value = 1692800
$ws.Range("I83").Value = 1650200
$ws.Range("J83").Value = "NA"

$ws.Range("D89").Value = 4333200
$ws.Range("E89").Value = 4676500
$ws.Range("F89").Value = 6287000
$ws.Range("G89").Value = 4945000
$ws.Range("H89").Value = 5971900
$ws.Range("I89").Value = 3996400
$ws.Range("J89").Value = 3109100

$ws.Range("D91").Value = -3828700
$ws.Range("E91").Value = -2989700
$ws.Range("F91").Value = -3590200
$ws.Range("G91").Value = -2947200
$ws.Range("H91").Value = -2625600
$ws.Range("I91").Value = -2296600
$ws.Range("J91").Value = -2474500

$ws.Range("D94").Value = -4711800
$ws.Range("E94").Value = -2619400
$ws.Range("F94").Value = -4899800
$ws.Range("G94").Value = -2889700
$ws.Range("H94").Value = -2827800
$ws.Range("I94").Value = -2272100
$ws.Range("J94").Value = "NA"

$ws.Range("D96").Value = -646000
$ws.Range("E96").Value = -1069500
$ws.Range("F96").Value = -821500
$ws.Range("G96").Value = -1994400
$ws.Range("H96").Value = -1668400
$ws.Range("I96").Value = -1431400
$ws.Range("J96").Value = -1625400

$ws.Range("D100").Value = 595600
$ws.Range("E100").Value = -2018700
$ws.Range("F100").Value = -2098600
$ws.Range("G100").Value = -1616800
$ws.Range("H100").Value = -3300700
$ws.Range("I100").Value = -1456900
$ws.Range("J100").Value = "NA"

$ws.Range("E101").Value = 10800
$ws.Range("F101").Value = 4900
$ws.Range("G101").Value = -8700
$ws.Range("H101").Value = -16100
$ws.Range("I101").Value = 22400
$ws.Range("J101").Value = "NA"

$ws.Range("D102").Value = 218400
$ws.Range("E102").Value = 49300
$ws.Range("F102").Value = -706400
$ws.Range("G102").Value = 429900
$ws.Range("H102").Value = -172700
$ws.Range("I102").Value = 289800
$ws.Range("J102").Value = -129700
